# The commit simplifies the document's style <w:docDefaults> block
# (in word/styles.xml) down to only the handful of properties that are
# not equal to Word's own built-in defaults, i.e. it strips the
# redundant/explicit-default run and paragraph properties that were
# being emitted for every document's default formatting.
#
# There is no dedicated Word object-model surface for docDefaults, so we
# reach it the same way Word's own COM automation would for "raw OOXML"
# edits: via Document.WordOpenXML, which round-trips the full package
# (including styles.xml) as a single XML string. We locate the
# <w:docDefaults>...</w:docDefaults> element and replace it wholesale
# with the trimmed-down version, then write the string back.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$startTag = "<w:docDefaults>"
$endTag = "</w:docDefaults>"

$startIdx = $xml.IndexOf($startTag)
$endIdx = $xml.IndexOf($endTag) + $endTag.Length

$newDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newXml = $xml.Substring(0, $startIdx) + $newDocDefaults + $xml.Substring($endIdx)

$d.WordOpenXML = $newXml
